# Fri, Mar 20, 2020  6:06:21 PM
#
# 1) The table on slide 6 gets a new (built-in) table style applied via the
#    Table Design gallery.
# 2) The presentation's colour theme is swapped back to the default
#    "Office" colour scheme (the deck keeps its Integral-derived fonts /
#    fill & effect formats - only the colour swatches change).

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 6 --------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B663A699-D182-41DF-867E-8E683740E37F}")
    }
}

# --- 2) Swap the theme colours back to the standard Office palette ---------
function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToRgbVal $officeColors[$i - 1]
}
